$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "RND_" prefix to the raw-data generator expressions in columns
# B (Seats), D (Frequency), E (Cov. Sum), F (Value) for rows 2-4.
foreach ($row in 2..4) {
    foreach ($col in @("B", "D", "E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value = "RND_" + $cell.Value2
    }
}
